$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorganizacion: renombrar la hoja "Datos" a "calidad_animal"
$ws.Name = "calidad_animal"

# Quitar el ancho de columna fijo (20) que tenian A:C, volviendo al
# comportamiento por defecto de la hoja (limpieza de formato antiguo)
$ws.Columns("A:C").Delete()
$ws.Columns("A:C").Insert()

# Encabezados en minusculas y sin acentos (nuevo esquema de carga)
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "comentario"

# Quitar el estilo (negrita blanca sobre relleno azul, centrado) que
# tenia la fila de encabezado
$ws.Range("A1:C1").ClearFormats()
